# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the financial-documents table) gets a new
#    table style applied (tableStyleId swap).
# 2) The deck's theme colour scheme is changed from the "Integral" /
#    "Red Violet" palette to the stock "Office" palette (the theme
#    that drives the slide master / all slides).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{378B8DA1-6F0E-4E38-87C1-385B0BB77C90}")

# --- 2. Swap the theme colour scheme over to "Office" ---------------------
function Set-ThemeColor {
    param($colorScheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

Set-ThemeColor $colorScheme 1  "000000"
Set-ThemeColor $colorScheme 2  "FFFFFF"
Set-ThemeColor $colorScheme 3  "44546A"
Set-ThemeColor $colorScheme 4  "E7E6E6"
Set-ThemeColor $colorScheme 5  "5B9BD5"
Set-ThemeColor $colorScheme 6  "ED7D31"
Set-ThemeColor $colorScheme 7  "A5A5A5"
Set-ThemeColor $colorScheme 8  "FFC000"
Set-ThemeColor $colorScheme 9  "4472C4"
Set-ThemeColor $colorScheme 10 "70AD47"
Set-ThemeColor $colorScheme 11 "0563C1"
Set-ThemeColor $colorScheme 12 "954F72"
